$d = $word.ActiveDocument

$replacements = @(
    @("127÷7=", "780÷2="),
    @("437÷7=", "477÷2="),
    @("335÷9=", "415÷8="),
    @("808÷7=", "271÷2="),
    @("734÷6=", "529÷2="),
    @("469÷4=", "201÷9="),
    @("784÷5=", "310÷7="),
    @("208÷5=", "911÷8="),
    @("834÷7=", "413÷8="),
    @("697÷7=", "947÷3="),
    @("897÷4=", "848÷9="),
    @("767÷4=", "682÷4="),
    @("113÷9=", "863÷7="),
    @("380÷8=", "970÷5="),
    @("498÷7=", "668÷6="),
    @("888÷6=", "168÷3="),
    @("489÷2=", "832÷6="),
    @("302÷4=", "727÷6="),
    @("202÷4=", "149÷6="),
    @("924÷5=", "337÷4="),
    @("120÷9=", "384÷5="),
    @("881÷6=", "912÷2="),
    @("893÷2=", "530÷7="),
    @("677÷8=", "907÷3="),
    @("498÷6=", "855÷7=")
)

foreach ($pair in $replacements) {
    $old = $pair[0]
    $new = $pair[1]
    $range = $d.Content
    $range.Find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, $new, 2)
}
